$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.62"
$ws.Range("E2").Value = "'0.51%"
$ws.Range("D3").Value = "'29.73"
$ws.Range("E3").Value = "'10.27%"
$ws.Range("D4").Value = "'5.174"
$ws.Range("E4").Value = "'1.73%"
$ws.Range("D5").Value = "'0.05719"
$ws.Range("E5").Value = "'0.40%"
$ws.Range("D6").Value = "'6.612"
$ws.Range("E6").Value = "'1.98%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.8568"
$ws.Range("E7").Value = "'4.40%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'0.8694"
$ws.Range("E8").Value = "'-3.52%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1367"
$ws.Range("E9").Value = "'2.69%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07076"
$ws.Range("E10").Value = "'2.43%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.02895"
$ws.Range("E11").Value = "'2.74%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09388"
$ws.Range("E12").Value = "'-0.09%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001514"
$ws.Range("E13").Value = "'0.31%"
$ws.Range("B14").Value = "CoinExToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D14").Value = "'0.04168"
$ws.Range("E14").Value = "'2.06%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006020"
$ws.Range("E15").Value = "'0.39%"
$ws.Range("D16").Value = "'0.005978"
$ws.Range("E16").Value = "'-1.14%"
$ws.Range("E17").Value = "'5,071.79%"
$ws.Range("D18").Value = "'3.488"
$ws.Range("E18").Value = "'-0.54%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.084"
$ws.Range("E19").Value = "'2.55%"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "'2.272"
$ws.Range("E20").Value = "'1.87%"
$ws.Range("E21").Value = "'-0.11%"
$ws.Range("D22").Value = "'0.03367"
$ws.Range("E22").Value = "'6.79%"
$ws.Range("E23").Value = "'2.10%"
$ws.Range("D24").Value = "'3.468"
$ws.Range("E24").Value = "'-2.36%"
$ws.Range("E25").Value = "'0.55%"
$ws.Range("D26").Value = "'0.005025"
$ws.Range("E26").Value = "'26.57%"
$ws.Range("D27").Value = "'0.001223"
$ws.Range("E27").Value = "'0.47%"
$ws.Range("E28").Value = "'23.53%"
$ws.Range("D40").Value = "'0.03745"
$ws.Range("E40").Value = "'1.25%"
$ws.Range("D41").Value = "'0.005770"
$ws.Range("E41").Value = "'68.04%"
$ws.Range("E42").Value = "'1.37%"
$ws.Range("E43").Value = "'-16.63%"
$ws.Range("D44").Value = "'0.008312"
$ws.Range("E44").Value = "'-11.56%"
$ws.Range("D45").Value = "'0.00005210"
$ws.Range("E45").Value = "'0.01%"
$ws.Range("E46").Value = "'0.06%"
$ws.Range("E47").Value = "'-51.64%"
$ws.Range("E48").Value = "'0.90%"
$ws.Range("E49").Value = "'0.06%"
$ws.Range("E50").Value = "'0.06%"
